# =========================================================================
# Restructure PlayerPerformance_4558.xlsx:
#   1. Add a new "Player Info" sheet before "ODI Batting"
#   2. Update "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, URLs -> bare codes,
#      drop the two blank INNING_NUMBER placeholder cells
#   3. Add a new "ODI Batting Extra" sheet after "ODI Batting"
# =========================================================================

$wb = $excel.ActiveWorkbook

# --- 1. Insert "Player Info" sheet before "ODI Batting" ---
$odi = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($odi)
$playerInfo.Name = "Player Info"

$playerInfo = $wb.Worksheets.Item("Player Info")
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
$playerInfo.Range("A1:D1").Font.Bold = $true
$playerInfo.Range("A1:D1").HorizontalAlignment = -4108
$playerInfo.Range("A1:D1").VerticalAlignment = -4160
$playerInfo.Range("A1:D1").Borders.LineStyle = 1

$playerInfo.Range("A2").Value = "'4558"
$playerInfo.Range("B2").Value = 'Peter Stephen Patrick Handscomb'
$playerInfo.Range("C2").Value = 'Right Handed'
$playerInfo.Range("D2").Value = 'Does Not Bowl | Unknown'

# --- 2. Update "ODI Batting": rename MATCH_CARD_LINK column + simplify D values ---
$odi = $wb.Worksheets.Item("ODI Batting")
$odi.Range("D1").Value = "MATCH_CODE"

# Replace full scorecard URLs with bare numeric match codes (stored as text)
$odi.Range("D2").Value = "'3975"
$odi.Range("D3").Value = "'3977"
$odi.Range("D4").Value = "'3981"
$odi.Range("D5").Value = "'3984"
$odi.Range("D6").Value = "'3988"
$odi.Range("D7").Value = "'4071"
$odi.Range("D8").Value = "'4074"
$odi.Range("D9").Value = "'4076"
$odi.Range("D10").Value = "'4234"
$odi.Range("D11").Value = "'4235"
$odi.Range("D12").Value = "'4236"
$odi.Range("D13").Value = "'4258"
$odi.Range("D14").Value = "'4263"
$odi.Range("D15").Value = "'4266"
$odi.Range("D16").Value = "'4268"
$odi.Range("D17").Value = "'4270"
$odi.Range("D18").Value = "'4273"
$odi.Range("D19").Value = "'4274"
$odi.Range("D20").Value = "'4275"
$odi.Range("D21").Value = "'4276"
$odi.Range("D22").Value = "'4277"
$odi.Range("D23").Value = "'4354"

# Remove the two blank INNING_NUMBER placeholder cells (did-not-bat rows)
$odi.Range("B3").ClearContents()
$odi.Range("B19").ClearContents()

# --- 3. Append "ODI Batting Extra" sheet after "ODI Batting" ---
$odi = $wb.Worksheets.Item("ODI Batting")
$extra = $wb.Worksheets.Add($null, $odi)
$extra.Name = "ODI Batting Extra"

$extra = $wb.Worksheets.Item("ODI Batting Extra")
$extra.Range("A1").Value = 'MATCH_CODE'
$extra.Range("B1").Value = 'BATTING_POSITION'
$extra.Range("C1").Value = 'NUM_4'
$extra.Range("D1").Value = 'NUM_6'
$extra.Range("E1").Value = 'PERCENT_RUNS_OF_TOTAL'
$extra.Range("F1").Value = 'MAN_OF_MATCH'
$extra.Range("A1:F1").Font.Bold = $true
$extra.Range("A1:F1").HorizontalAlignment = -4108
$extra.Range("A1:F1").VerticalAlignment = -4160
$extra.Range("A1:F1").Borders.LineStyle = 1

# row 2: MATCH_CODE 3981
$extra.Range("A2").Value = "'3981"
$extra.Range("B2").Value = 6
$extra.Range("C2").Value = "'0"
$extra.Range("D2").Value = "'0"
$extra.Range("E2").Value = '0.27%'
$extra.Range("F2").Value = 'NO'

# row 3: MATCH_CODE 3984
$extra.Range("A3").Value = "'3984"
$extra.Range("B3").Value = 4
$extra.Range("C3").Value = "'1"
$extra.Range("D3").Value = "'0"
$extra.Range("E3").Value = '2.50%'
$extra.Range("F3").Value = 'NO'

# row 4: MATCH_CODE 3988
$extra.Range("A4").Value = "'3988"
$extra.Range("B4").Value = 3
$extra.Range("C4").Value = "'0"
$extra.Range("D4").Value = "'0"
$extra.Range("F4").Value = 'NO'

# row 5: MATCH_CODE 4071
$extra.Range("A5").Value = "'4071"
$extra.Range("B5").Value = 7
$extra.Range("C5").Value = "'0"
$extra.Range("D5").Value = "'0"
$extra.Range("E5").Value = '1.02%'
$extra.Range("F5").Value = 'NO'

# row 6: MATCH_CODE 4074
$extra.Range("A6").Value = "'4074"
$extra.Range("F6").Value = 'NO'

# row 7: MATCH_CODE 4076
$extra.Range("A7").Value = "'4076"
$extra.Range("B7").Value = 4
$extra.Range("C7").Value = "'0"
$extra.Range("D7").Value = "'0"
$extra.Range("E7").Value = '5.37%'
$extra.Range("F7").Value = 'NO'

# row 8: MATCH_CODE 4234
$extra.Range("A8").Value = "'4234"
$extra.Range("B8").Value = 5
$extra.Range("C8").Value = "'6"
$extra.Range("D8").Value = "'2"
$extra.Range("E8").Value = '25.35%'
$extra.Range("F8").Value = 'NO'

# row 9: MATCH_CODE 4235
$extra.Range("A9").Value = "'4235"
$extra.Range("F9").Value = 'NO'

# row 10: MATCH_CODE 4236
$extra.Range("A10").Value = "'4236"
$extra.Range("B10").Value = 5
$extra.Range("C10").Value = "'2"
$extra.Range("D10").Value = "'0"
$extra.Range("E10").Value = '25.22%'
$extra.Range("F10").Value = 'NO'

# row 11: MATCH_CODE 4258
$extra.Range("A11").Value = "'4258"
$extra.Range("F11").Value = 'NO'

# row 12: MATCH_CODE 4263
$extra.Range("A12").Value = "'4263"
$extra.Range("B12").Value = 4
$extra.Range("C12").Value = "'4"
$extra.Range("D12").Value = "'0"
$extra.Range("E12").Value = '19.83%'
$extra.Range("F12").Value = 'NO'

# row 13: MATCH_CODE 4266
$extra.Range("A13").Value = "'4266"
$extra.Range("B13").Value = 6
$extra.Range("C13").Value = "'0"
$extra.Range("D13").Value = "'0"
$extra.Range("F13").Value = 'NO'

# row 14: MATCH_CODE 4268
$extra.Range("A14").Value = "'4268"
$extra.Range("B14").Value = 4
$extra.Range("C14").Value = "'8"
$extra.Range("D14").Value = "'3"
$extra.Range("E14").Value = '32.59%'
$extra.Range("F14").Value = 'NO'

# row 15: MATCH_CODE 4270
$extra.Range("A15").Value = "'4270"
$extra.Range("B15").Value = 3
$extra.Range("C15").Value = "'4"
$extra.Range("D15").Value = "'0"
$extra.Range("E15").Value = '19.12%'
$extra.Range("F15").Value = 'NO'

# row 16: MATCH_CODE 4273
$extra.Range("A16").Value = "'4273"
$extra.Range("B16").Value = 4
$extra.Range("C16").Value = "'2"
$extra.Range("D16").Value = "'0"
$extra.Range("E16").Value = '10.68%'
$extra.Range("F16").Value = 'NO'

# row 17: MATCH_CODE 4274
$extra.Range("A17").Value = "'4274"
$extra.Range("B17").Value = 5
$extra.Range("F17").Value = 'NO'

# row 18: MATCH_CODE 4275
$extra.Range("A18").Value = "'4275"
$extra.Range("B18").Value = 4
$extra.Range("C18").Value = "'6"
$extra.Range("D18").Value = "'0"
$extra.Range("E18").Value = '17.67%'
$extra.Range("F18").Value = 'NO'

# row 19: MATCH_CODE 4276
$extra.Range("A19").Value = "'4276"
$extra.Range("F19").Value = 'NO'

# row 20: MATCH_CODE 4277
$extra.Range("A20").Value = "'4277"
$extra.Range("B20").Value = 6
$extra.Range("C20").Value = "'1"
$extra.Range("D20").Value = "'0"
$extra.Range("E20").Value = '2.45%'
$extra.Range("F20").Value = 'NO'

# row 21: MATCH_CODE 4354
$extra.Range("A21").Value = "'4354"
$extra.Range("B21").Value = 4
$extra.Range("C21").Value = "'0"
$extra.Range("D21").Value = "'0"
$extra.Range("E21").Value = '1.79%'
$extra.Range("F21").Value = 'NO'

